$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits on the "Kick-Off Meeting"
#    heading paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Turn the first of the two empty bold paragraphs that follow the
#    "Abstract" paragraph into a tab-stop paragraph (adds a tab stop + a
#    bold tab run), give the "11:09 a.m. - Meeting start." paragraph a
#    hanging indent and append the "All present: ..." roster of names to
#    it, and insert a new empty bold paragraph right after it.
# ---------------------------------------------------------------------------
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$tabParagraphXml = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2636"/></w:tabs><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r></w:p>'

$namesRunsXml = '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">All present: </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Kevin O' + [char]0x2019 + 'Hare, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Conor</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>McAleavey</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">, Jonnie </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Leathem</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">, Steven Kennedy, Chris </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>McClun</w:t></w:r><w:r><w:t>e</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">, Alan Whitten, Jack Ferguson, </w:t></w:r>' + `
  '<w:r><w:t>Adam Hale</w:t></w:r>'

$meetingStartParagraphXml = '<w:p><w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr>' + `
  '<w:r w:rsidRPr="00934310"><w:rPr><w:b/></w:rPr><w:t>11:09</w:t></w:r>' + `
  '<w:r w:rsidR="004F252D"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> a.m.</w:t></w:r>' + `
  '<w:r w:rsidRPr="00934310"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> ' + [char]0x2013 + '</w:t></w:r>' + `
  '<w:r w:rsidR="00934310" w:rsidRPr="00934310"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Meeting start. </w:t></w:r>' + `
  $namesRunsXml + `
  '</w:p>' + `
  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>'

# Find the two source paragraphs by their (still-plain) text before editing.
$p6 = $d.Paragraphs.Item(6)
$p7 = $d.Paragraphs.Item(7)

$p6.Range.InsertXML("<w:document $xmlNs>$tabParagraphXml</w:document>")
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertXML("<w:document $xmlNs>$meetingStartParagraphXml</w:document>")

# ---------------------------------------------------------------------------
# 3) Split the "We intend to meet ... meetings, especially as issues
#    arise." sentence and re-insert a _GoBack bookmark at the split point.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("real-world processes for ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $findRng.End
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 4) Bold the "11:37 a.m:" paragraph (paragraph mark + every run).
# ---------------------------------------------------------------------------
$timeRng = $d.Content
$timeRng.Find.Execute("11:37 ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$timeParagraph = $d.Paragraphs.Item($timeRng.Paragraphs.First.Index)
$timeParagraph.Range.Font.Bold = 1
